# Add "better support for importing capital calls & distributions":
# two new columns -- Generate Remittances (E) and Remittances Verified (F) --
# with Yes/No flags for each existing capital-call row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("E1").Value = "Generate Remittances"
$ws.Range("F1").Value = "Remittances Verified"

# Data rows
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "Yes"

$ws.Range("E3").Value = "Yes"
$ws.Range("F3").Value = "No"

$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "No"

# Best-fit the two new columns, same as the existing A:D columns
$ws.Columns.Item(5).ColumnWidth = 18.3
$ws.Columns.Item(6).ColumnWidth = 17.03

# Match the saved selection/active cell from the source workbook
$ws.Range("E5").Select()
